# Apply "Natmi following Dr Hou advice" data refresh to sheet1.
# Recomputed Clcf1 -> Lifr LR-pair table across all three clusters (ECs, FAPs, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Clcf1"
$ws.Cells.Item(2, 3).Value = "Lifr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.798783666666667
$ws.Cells.Item(2, 8).Value = 5.396351
$ws.Cells.Item(2, 9).Value = 0.2319744053785674
$ws.Cells.Item(2, 10).Value = 0.2319744053785674
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 14.96835266666667
$ws.Cells.Item(2, 14).Value = 44.905058
$ws.Cells.Item(2, 15).Value = 0.1240053612000741
$ws.Cells.Item(2, 16).Value = 0.1240053612000741
$ws.Cells.Item(2, 17).Value = 26.92482829370644
$ws.Cells.Item(2, 18).Value = 242.323454643358
$ws.Cells.Item(2, 19).Value = 0.02876606992814166
$ws.Cells.Item(2, 20).Value = 0.02876606992814166

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Clcf1"
$ws.Cells.Item(3, 3).Value = "Lifr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.798783666666667
$ws.Cells.Item(3, 8).Value = 5.396351
$ws.Cells.Item(3, 9).Value = 0.2319744053785674
$ws.Cells.Item(3, 10).Value = 0.2319744053785674
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 80.77474466666666
$ws.Cells.Item(3, 14).Value = 242.324234
$ws.Cells.Item(3, 15).Value = 0.6691786071115035
$ws.Cells.Item(3, 16).Value = 0.6691786071115035
$ws.Cells.Item(3, 17).Value = 145.2962913855704
$ws.Cells.Item(3, 18).Value = 1307.666622470134
$ws.Cells.Item(3, 19).Value = 0.155232309476749
$ws.Cells.Item(3, 20).Value = 0.155232309476749

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Clcf1"
$ws.Cells.Item(4, 3).Value = "Lifr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.798783666666667
$ws.Cells.Item(4, 8).Value = 5.396351
$ws.Cells.Item(4, 9).Value = 0.2319744053785674
$ws.Cells.Item(4, 10).Value = 0.2319744053785674
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 24.96420533333334
$ws.Cells.Item(4, 14).Value = 74.892616
$ws.Cells.Item(4, 15).Value = 0.2068160316884225
$ws.Cells.Item(4, 16).Value = 0.2068160316884225
$ws.Cells.Item(4, 17).Value = 44.9052048049129
$ws.Cells.Item(4, 18).Value = 404.146843244216
$ws.Cells.Item(4, 19).Value = 0.04797602597367676
$ws.Cells.Item(4, 20).Value = 0.04797602597367676

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Clcf1"
$ws.Cells.Item(5, 3).Value = "Lifr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.782700333333333
$ws.Cells.Item(5, 8).Value = 5.348101
$ws.Cells.Item(5, 9).Value = 0.2299002695301921
$ws.Cells.Item(5, 10).Value = 0.2299002695301921
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 14.96835266666667
$ws.Cells.Item(5, 14).Value = 44.905058
$ws.Cells.Item(5, 15).Value = 0.1240053612000741
$ws.Cells.Item(5, 16).Value = 0.1240053612000741
$ws.Cells.Item(5, 17).Value = 26.68408728831756
$ws.Cells.Item(5, 18).Value = 240.156785594858
$ws.Cells.Item(5, 19).Value = 0.02850886596308586
$ws.Cells.Item(5, 20).Value = 0.02850886596308586

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Clcf1"
$ws.Cells.Item(6, 3).Value = "Lifr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.782700333333333
$ws.Cells.Item(6, 8).Value = 5.348101
$ws.Cells.Item(6, 9).Value = 0.2299002695301921
$ws.Cells.Item(6, 10).Value = 0.2299002695301921
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 80.77474466666666
$ws.Cells.Item(6, 14).Value = 242.324234
$ws.Cells.Item(6, 15).Value = 0.6691786071115035
$ws.Cells.Item(6, 16).Value = 0.6691786071115035
$ws.Cells.Item(6, 17).Value = 143.9971642421816
$ws.Cells.Item(6, 18).Value = 1295.974478179634
$ws.Cells.Item(6, 19).Value = 0.1538443421387732
$ws.Cells.Item(6, 20).Value = 0.1538443421387732

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Clcf1"
$ws.Cells.Item(7, 3).Value = "Lifr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.782700333333333
$ws.Cells.Item(7, 8).Value = 5.348101
$ws.Cells.Item(7, 9).Value = 0.2299002695301921
$ws.Cells.Item(7, 10).Value = 0.2299002695301921
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 24.96420533333334
$ws.Cells.Item(7, 14).Value = 74.892616
$ws.Cells.Item(7, 15).Value = 0.2068160316884225
$ws.Cells.Item(7, 16).Value = 0.2068160316884225
$ws.Cells.Item(7, 17).Value = 44.50369716913512
$ws.Cells.Item(7, 18).Value = 400.533274522216
$ws.Cells.Item(7, 19).Value = 0.04754706142833308
$ws.Cells.Item(7, 20).Value = 0.04754706142833308

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Clcf1"
$ws.Cells.Item(8, 3).Value = "Lifr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.172749333333334
$ws.Cells.Item(8, 8).Value = 12.518248
$ws.Cells.Item(8, 9).Value = 0.5381253250912406
$ws.Cells.Item(8, 10).Value = 0.5381253250912404
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 14.96835266666667
$ws.Cells.Item(8, 14).Value = 44.905058
$ws.Cells.Item(8, 15).Value = 0.1240053612000741
$ws.Cells.Item(8, 16).Value = 0.1240053612000741
$ws.Cells.Item(8, 17).Value = 62.45918361093156
$ws.Cells.Item(8, 18).Value = 562.132652498384
$ws.Cells.Item(8, 19).Value = 0.06673042530884658
$ws.Cells.Item(8, 20).Value = 0.06673042530884657

# Row 9: sCs -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Clcf1"
$ws.Cells.Item(9, 3).Value = "Lifr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.172749333333334
$ws.Cells.Item(9, 8).Value = 12.518248
$ws.Cells.Item(9, 9).Value = 0.5381253250912406
$ws.Cells.Item(9, 10).Value = 0.5381253250912404
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 80.77474466666666
$ws.Cells.Item(9, 14).Value = 242.324234
$ws.Cells.Item(9, 15).Value = 0.6691786071115035
$ws.Cells.Item(9, 16).Value = 0.6691786071115035
$ws.Cells.Item(9, 17).Value = 337.0527619580036
$ws.Cells.Item(9, 18).Value = 3033.474857622032
$ws.Cells.Item(9, 19).Value = 0.3601019554959813
$ws.Cells.Item(9, 20).Value = 0.3601019554959813

# Row 10: sCs -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Clcf1"
$ws.Cells.Item(10, 3).Value = "Lifr"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.172749333333334
$ws.Cells.Item(10, 8).Value = 12.518248
$ws.Cells.Item(10, 9).Value = 0.5381253250912406
$ws.Cells.Item(10, 10).Value = 0.5381253250912404
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 24.96420533333334
$ws.Cells.Item(10, 14).Value = 74.892616
$ws.Cells.Item(10, 15).Value = 0.2068160316884225
$ws.Cells.Item(10, 16).Value = 0.2068160316884225
$ws.Cells.Item(10, 17).Value = 104.1693711618631
$ws.Cells.Item(10, 18).Value = 937.5243404567682
$ws.Cells.Item(10, 19).Value = 0.1112929442864126
$ws.Cells.Item(10, 20).Value = 0.1112929442864126

